# WCP.docx edit: replace the 1x1 placeholder picture (the "submission
# flowchart" image stand-in) with a plain hyperlink run pointing at the
# hosted flowchart image, using the same visible text as the URL itself.
#
# Everything else in the document (headings, body text, the table, the
# three pre-existing external hyperlinks, and all of the bookmarks) is
# left untouched.

$d = $word.ActiveDocument

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg"

# The picture is the document's only inline shape - grab its Range before
# removing it so the Range can be reused as the insertion point for the
# new hyperlink run.
$pic = $d.InlineShapes.Item(1)
$picRange = $pic.Range
$pic.Delete()

# Insert the hyperlink in place of the (now deleted) picture. Leaving
# TextToDisplay unset makes Word display the address itself, matching
# "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg".
$d.Hyperlinks.Add($picRange, $url) | Out-Null
